# Re-run of the parser: refresh rows 2-3 with corrected values and append rows 4-7
# for the two new source records (41645857, 42377157).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that no longer carry data in rows 2-3 ---
$ws.Range("J2").Value = ""
$ws.Range("Q2").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("Q3").Value = ""

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 41645549
$ws.Range("D2").Value = "https://elibrary.ru/item.asp?id=41645549"
$ws.Range("E2").Value = "статья в журнале"
$ws.Range("F2").Value = "материалы конференции"
$ws.Range("G2").Value = "Lecture Notes in Etetete"
$ws.Range("H2").Value = "1876-1100"
$ws.Range("I2").Value = "1876-1119"
$ws.Range("K2").Value = "no"
$ws.Range("M2").Value = "no"
$ws.Range("N2").Value = "yes"
$ws.Range("P2").Value = 2020
$ws.Range("S2").Value = 554
$ws.Range("T2").Value = 154
$ws.Range("U2").Value = 161
$ws.Range("V2").Value = "EN"
$ws.Range("W2").Value = "Geometrical computational method to locate hypocenter by signal readings from a three receivers"
$ws.Range("X2").Value = "10.1007/978-3-030-14907-9_16"
$ws.Range("Y2").Value = "YMFJKJ"
$ws.Range("Z2").Value = 373100
$ws.Range("AA2").Value = "yes"
$ws.Range("AB2").Value = "yes"
$ws.Range("AC2").Value = 2
$ws.Range("AD2").Value = "excel"

# --- Row 3 ---
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 41645549
$ws.Range("D3").Value = "https://elibrary.ru/item.asp?id=41645549"
$ws.Range("E3").Value = "статья в журнале"
$ws.Range("F3").Value = "материалы конференции"
$ws.Range("G3").Value = "Lecture Notes in Electrical Engineering"
$ws.Range("H3").Value = "1876-1100"
$ws.Range("I3").Value = "1876-1119"
$ws.Range("K3").Value = "no"
$ws.Range("M3").Value = "no"
$ws.Range("N3").Value = "yes"
$ws.Range("P3").Value = 2020
$ws.Range("S3").Value = 554
$ws.Range("T3").Value = 154
$ws.Range("U3").Value = 161
$ws.Range("V3").Value = "EN"
$ws.Range("W3").Value = "Geometrical computational method to locate hypocenter by signal readings from a three receivers"
$ws.Range("X3").Value = "10.1007/978-3-030-14907-9_16"
$ws.Range("Y3").Value = "YMFJKJ"
$ws.Range("Z3").Value = 373100
$ws.Range("AA3").Value = "yes"
$ws.Range("AB3").Value = "yes"
$ws.Range("AC3").Value = 2
$ws.Range("AD3").Value = "sql"

# --- Row 4 ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 41645857
$ws.Range("D4").Value = "https://elibrary.ru/item.asp?id=41645857"
$ws.Range("E4").Value = "статья в журнале"
$ws.Range("F4").Value = "материалы конференции"
$ws.Range("G4").Value = "Lecture fsfsfsfs"
$ws.Range("H4").Value = "1876-1100"
$ws.Range("I4").Value = "1876-1119"
$ws.Range("K4").Value = "no"
$ws.Range("M4").Value = "no"
$ws.Range("N4").Value = "yes"
$ws.Range("P4").Value = 2020
$ws.Range("S4").Value = 554
$ws.Range("T4").Value = 185
$ws.Range("U4").Value = 194
$ws.Range("V4").Value = "EN"
$ws.Range("W4").Value = "An investigation on signal comparison by measuring of numerical strings similarity"
$ws.Range("X4").Value = "10.1007/978-3-030-14907-9_19"
$ws.Range("Y4").Value = "WOVYTY"
$ws.Range("Z4").Value = 290000
$ws.Range("AA4").Value = "yes"
$ws.Range("AB4").Value = "yes"
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = "excel"

# --- Row 5 ---
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 414
$ws.Range("C5").Value = 41645857
$ws.Range("D5").Value = "https://elibrary.ru/item.asp?id=41645857"
$ws.Range("E5").Value = "статья в журнале"
$ws.Range("F5").Value = "материалы конференции"
$ws.Range("G5").Value = "Lecture Notes in Electrical Engineering"
$ws.Range("H5").Value = "1876-1100"
$ws.Range("I5").Value = "1876-1119"
$ws.Range("K5").Value = "no"
$ws.Range("M5").Value = "no"
$ws.Range("N5").Value = "yes"
$ws.Range("P5").Value = 2020
$ws.Range("S5").Value = 554
$ws.Range("T5").Value = 185
$ws.Range("U5").Value = 194
$ws.Range("V5").Value = "EN"
$ws.Range("W5").Value = "An investigation on signal comparison by measuring of numerical strings similarity"
$ws.Range("X5").Value = "10.1007/978-3-030-14907-9_19"
$ws.Range("Y5").Value = "WOVYTY"
$ws.Range("Z5").Value = 290000
$ws.Range("AA5").Value = "yes"
$ws.Range("AB5").Value = "yes"
$ws.Range("AC5").Value = 6
$ws.Range("AD5").Value = "sql"

# --- Row 6 ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = 42377157
$ws.Range("D6").Value = "https://elibrary.ru/item.asp?id=42377157"
$ws.Range("E6").Value = "статья в журнале"
$ws.Range("F6").Value = "научная статья"
$ws.Range("G6").Value = "Нефтепромысловое делоshhshsh"
$ws.Range("H6").Value = "0207-2351"
$ws.Range("J6").Value = "Российский государственный университет нефти и газа (национальный исследовательский университет) им. И.М. Губкина"
$ws.Range("K6").Value = "yes"
$ws.Range("M6").Value = "no"
$ws.Range("N6").Value = "no"
$ws.Range("P6").Value = 2020
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 614
$ws.Range("T6").Value = 32
$ws.Range("U6").Value = 37
$ws.Range("V6").Value = "RU"
$ws.Range("W6").Value = "Тестирование полимерно-гелевых систем `"Темпоскрин-Плюс ВПП`" и `"Темпоскрин-Люкс`" в условиях применения высокоминерализованных агентов закачки при высоких температурах пласта"
$ws.Range("X6").Value = "10.30713/0207-2351-2020-2(614)-32-37"
$ws.Range("Y6").Value = "UBTMHM"
$ws.Range("Z6").Value = 524700
$ws.Range("AA6").Value = "yes"
$ws.Range("AB6").Value = "no"
$ws.Range("AC6").Value = 22
$ws.Range("AD6").Value = "excel"

# --- Row 7 ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 377
$ws.Range("C7").Value = 42377157
$ws.Range("D7").Value = "https://elibrary.ru/item.asp?id=42377157"
$ws.Range("E7").Value = "статья в журнале"
$ws.Range("F7").Value = "научная статья"
$ws.Range("G7").Value = "Нефтепромысловое дело"
$ws.Range("H7").Value = "0207-2351"
$ws.Range("J7").Value = "Российский государственный университет нефти и газа (национальный исследовательский университет) им. И.М. Губкина"
$ws.Range("K7").Value = "yes"
$ws.Range("M7").Value = "no"
$ws.Range("N7").Value = "no"
$ws.Range("P7").Value = 2020
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 614
$ws.Range("T7").Value = 32
$ws.Range("U7").Value = 37
$ws.Range("V7").Value = "RU"
$ws.Range("W7").Value = "Тестирование полимерно-гелевых систем `"Темпоскрин-Плюс ВПП`" и `"Темпоскрин-Люкс`" в условиях применения высокоминерализованных агентов закачки при высоких температурах пласта"
$ws.Range("X7").Value = "10.30713/0207-2351-2020-2(614)-32-37"
$ws.Range("Y7").Value = "UBTMHM"
$ws.Range("Z7").Value = 524700
$ws.Range("AA7").Value = "yes"
$ws.Range("AB7").Value = "no"
$ws.Range("AC7").Value = 22
$ws.Range("AD7").Value = "sql"

# --- Propagate the header-ish row style (bold, centered, thin border) used in column A ---
$ws.Range("A2").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rebuild the linkurl hyperlinks cleanly (avoids leaving stale/duplicate link refs
#     behind when re-pointing D2/D3 at their corrected URLs) ---
$ws.Range("D2").Copy()
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "https://elibrary.ru/item.asp?id=41645549") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://elibrary.ru/item.asp?id=41645549") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://elibrary.ru/item.asp?id=41645857") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://elibrary.ru/item.asp?id=41645857") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://elibrary.ru/item.asp?id=42377157") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://elibrary.ru/item.asp?id=42377157") | Out-Null

# Hyperlinks.Add mints its own (near-identical) style; restore the original hyperlink
# cell style (captured above) across every linkurl cell.
$ws.Range("D2:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done"
